$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Color constants (Interior.Color uses BGR ordering, i.e. 0xBBGGRR)
$green  = 5296274   # RGB(0x92,0xD0,0x50) -> BGR 0x50D092 (the existing "Done" row highlight)
$yellow = 65535      # RGB(0xFF,0xFF,0x00) -> BGR 0x00FFFF (the existing "Active" row highlight)

# Row 8 "Add Modal" is now finished
$ws.Range("C8").Value = "Done"

# Row 9 "Products Backend" is now finished (was the active task, now done)
$ws.Range("A9:C9").Interior.Color = $green
$ws.Range("C9").Value = "Done"

# Row 10 "Add Redux" is now finished
$ws.Range("A10:C10").Interior.Color = $green
$ws.Range("C10").Value = "Done"

# Row 11 "Add Redux To Products" is now finished
$ws.Range("A11:C11").Interior.Color = $green
$ws.Range("C11").Value = "Done"

# Row 12 "Add Redux To Filter" is now the active task
$ws.Range("A12:D12").Interior.Color = $yellow
$ws.Range("C12").Value = "Active"

# Move the selection to the new active task
$ws.Range("A12").Select()
